# Applies commit "Update gh-pages to output generated at 456a3b4"
# Updates column F ("想去人数") values across all 4 sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1563
$ws.Range("F3").Value = 941
$ws.Range("F4").Value = 930
$ws.Range("F5").Value = 529
$ws.Range("F6").Value = 7987
$ws.Range("F8").Value = 39
$ws.Range("F9").Value = 1941
$ws.Range("F10").Value = 5764
$ws.Range("F11").Value = 581
$ws.Range("F14").Value = 8135
$ws.Range("F15").Value = 9480
$ws.Range("F16").Value = 1160
$ws.Range("F17").Value = 944
$ws.Range("F18").Value = 4580
$ws.Range("F19").Value = 708
$ws.Range("F20").Value = 284
$ws.Range("F22").Value = 295
$ws.Range("F24").Value = 1225
$ws.Range("F25").Value = 136
$ws.Range("F26").Value = 1731
$ws.Range("F27").Value = 766
$ws.Range("F28").Value = 995
$ws.Range("F29").Value = 32
$ws.Range("F30").Value = 1922
$ws.Range("F32").Value = 491
$ws.Range("F33").Value = 2389
$ws.Range("F35").Value = 124
$ws.Range("F36").Value = 1520
$ws.Range("F38").Value = 1333
$ws.Range("F39").Value = 11
$ws.Range("F40").Value = 823
$ws.Range("F41").Value = 532
$ws.Range("F42").Value = 205
$ws.Range("F44").Value = 443
$ws.Range("F45").Value = 527
$ws.Range("F46").Value = 18
$ws.Range("F47").Value = 865
$ws.Range("F49").Value = 4123

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 25
$ws.Range("F9").Value = 25
$ws.Range("F15").Value = 51

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5423

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1563
$ws.Range("F3").Value = 941
$ws.Range("F4").Value = 930
$ws.Range("F5").Value = 529
$ws.Range("F7").Value = 39
$ws.Range("F9").Value = 25
$ws.Range("F11").Value = 5764
$ws.Range("F12").Value = 581
$ws.Range("F13").Value = 8135
$ws.Range("F14").Value = 9480
$ws.Range("F16").Value = 1160
$ws.Range("F17").Value = 944
$ws.Range("F18").Value = 708
$ws.Range("F19").Value = 284
$ws.Range("F21").Value = 295
$ws.Range("F23").Value = 51
$ws.Range("F24").Value = 1225
$ws.Range("F25").Value = 136
$ws.Range("F26").Value = 766
$ws.Range("F27").Value = 995
$ws.Range("F28").Value = 32
$ws.Range("F29").Value = 1922
$ws.Range("F31").Value = 491
$ws.Range("F32").Value = 2389
$ws.Range("F40").Value = 532
$ws.Range("F42").Value = 205
$ws.Range("F44").Value = 443
$ws.Range("F45").Value = 527
$ws.Range("F46").Value = 18
$ws.Range("F47").Value = 865
$ws.Range("F49").Value = 4123
